# Add new time-log entries (move and turn) to the "Сессии" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Сессии")

$ws.Cells.Item(14, 3).Value = 0.375
$ws.Cells.Item(14, 4).Value = "ребёнок"

$ws.Cells.Item(15, 3).Value = 0.37986111111111115

$ws.Cells.Item(16, 3).Value = 0.38125000000000003
$ws.Cells.Item(16, 4).Value = "ScrollLock"

$ws.Cells.Item(17, 3).Value = 0.38194444444444442

$ws.Cells.Item(18, 3).Value = 0.3840277777777778

$ws.Range("C20").Select() | Out-Null
